# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (2..24), a new "latest" error value is inserted at the
# front of the series (column B). Every existing value in that row shifts
# one column to the right (B->C, C->D, ... J->K). Any value that would have
# shifted past column K (the 10th data column) falls off the end.
#
# Row 2 is a special case: it already has a value in every column B..K, and
# no new value is being inserted there - the diff only drops the trailing K2
# value (the oldest point in that particular series has aged out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert into column B for each row (the newest data point).
$newValues = @{
    3  = 2.174397541324862
    4  = 8.587172912142332
    5  = -8.713442361886736
    6  = -0.58442257821662
    7  = 1.454533757567239
    8  = -1.777394389465022
    9  = -1.722070219091221
    10 = 0.3648791949059138
    11 = -0.2352699264540507
    12 = -0.05148746350304451
    13 = -0.1333319740152609
    14 = 1.614150253737389
    15 = 0.5701030647716323
    16 = 0.2202779152847414
    17 = 0.5040960054549828
    18 = 0.420735823599318
    19 = -0.1252583916527783
    20 = 0.08824118641116785
    21 = -0.1133200159455487
    22 = 0.1743923273248104
    23 = -0.4559694969238889
    24 = 0.1808172637304477
}

$lastDataCol = 11   # column K - the widest a row is ever allowed to grow

# Row 2: no insertion, the trailing value in column K simply drops off.
$ws.Range("K2").ClearContents()

# Rows 3..24: shift existing B..(last used column) values right by one
# column, then drop the new value into column B.
for ($r = 3; $r -le 24; $r++) {

    # Find the last used column in this row before editing.
    $lastCol = 1
    for ($c = 2; $c -le $lastDataCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -ne $null) {
            $lastCol = $c
        }
    }

    # Read the existing row values (columns 2..$lastCol) first.
    $rowValues = @{}
    for ($c = 2; $c -le $lastCol; $c++) {
        $rowValues[$c] = $ws.Cells.Item($r, $c).Value2
    }

    # Write them back shifted one column to the right, dropping anything
    # that would spill past column K.
    for ($c = $lastCol; $c -ge 2; $c--) {
        $destCol = $c + 1
        if ($destCol -le $lastDataCol) {
            $ws.Cells.Item($r, $destCol).Value = $rowValues[$c]
        }
    }

    # Insert the new leading value.
    $ws.Cells.Item($r, 2).Value = $newValues[$r]
}
